# Applies the crypto price/volume updates described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "59.779.58"
$ws.Range("E2").Value = "  -0.17%  "

# Row 3
$ws.Range("D3").Value = "2.377.01"
$ws.Range("E3").Value = "  -1.53%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").Value = "'555.73"
$ws.Range("E5").Value = "  +0.74%  "

# Row 6
$ws.Range("E6").Value = "  -2.85%  "

# Row 7
$ws.Range("E7").Value = "  +0.14%  "

# Row 8
$ws.Range("D8").Value = "'0.586"

# Row 9
$ws.Range("E9").Value = "  -0.45%  "

# Row 10
$ws.Range("D10").Value = "'5.62"
$ws.Range("E10").Value = "  -1.35%  "

# Row 11
$ws.Range("E11").Value = "  +1.10%  "

# Row 12
$ws.Range("E12").Value = "  -2.95%  "

# Row 13
$ws.Range("D13").Value = "'24.39"
$ws.Range("E13").Value = "  -4.86%  "

# Row 14
$ws.Range("D14").Value = "2.806.99"
$ws.Range("E14").Value = "  -1.28%  "

# Row 15
$ws.Range("D15").Value = "59.756.97"
$ws.Range("E15").Value = "  -0.09%  "

# Row 16
$ws.Range("D16").Value = "'0.0000136"
$ws.Range("E16").Value = "  -0.69%  "

# Row 17
$ws.Range("D17").Value = "2.383.57"
$ws.Range("E17").Value = "  -4.77%  "

# Row 18
$ws.Range("D18").Value = "'11.10"
$ws.Range("E18").Value = "  -1.90%  "

# Row 19
$ws.Range("D19").Value = "'4.47"
$ws.Range("E19").Value = "  +1.33%  "

# Row 20
$ws.Range("D20").Value = "'320.55"
$ws.Range("E20").Value = "  -2.55%  "

# Row 21
$ws.Range("D21").Value = "'6.69"
$ws.Range("E21").Value = "  +0.59%  "

# Row 22
$ws.Range("E22").Value = "  -0.02%  "

# Row 23
$ws.Range("D23").Value = "'64.12"
$ws.Range("E23").Value = "  -3.54%  "

# Row 24
$ws.Range("D24").Value = "'0.173"
$ws.Range("E24").Value = "  +0.89%  "

# Row 25
$ws.Range("E25").Value = "  +0.16%  "

# Row 26
$ws.Range("D26").Value = "'8.37"
$ws.Range("E26").Value = "  -3.18%  "

# Row 27
$ws.Range("D27").Value = "'1.37"
$ws.Range("E27").Value = "  -0.01%  "

# Row 28
$ws.Range("D28").Value = "'1.78"
$ws.Range("E28").Value = "  +0.84%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0757"
$ws.Range("E29").Value = "  -2.31%  "

# Row 30
$ws.Range("D30").Value = "'169.52"
$ws.Range("E30").Value = "  +0.85%  "

# Row 31
$ws.Range("D31").Value = "'6.03"
$ws.Range("E31").Value = "  -1.20%  "

# Row 32
$ws.Range("D32").Value = "'1.08"
$ws.Range("E32").Value = "  +6.52%  "

# Row 33
$ws.Range("D33").Value = "'0.394"
$ws.Range("E33").Value = "  -3.47%  "

# Row 34
$ws.Range("D34").Value = "'18.15"
$ws.Range("E34").Value = "  -2.55%  "

# Row 35
$ws.Range("E35").Value = "  +0.00%  "

# Row 36
$ws.Range("E36").Value = "  +1.18%  "

# Row 37
$ws.Range("E37").Value = "  +0.02%  "

# Row 38
$ws.Range("E38").Value = "  -2.13%  "

# Row 39
$ws.Range("D39").Value = "'318.04"
$ws.Range("E39").Value = "  +1.53%  "

# Row 40
$ws.Range("E40").Value = "  -1.73%  "

# Row 41
$ws.Range("D41").Value = "'38.59"
$ws.Range("E41").Value = "  -2.41%  "

# Row 42
$ws.Range("D42").Value = "'145.77"
$ws.Range("E42").Value = "  +5.13%  "

# Row 43
$ws.Range("E43").Value = "  -4.41%  "

# Row 44
$ws.Range("D44").Value = "'0.0969"
$ws.Range("E44").Value = "  -0.01%  "

# Row 45
$ws.Range("D45").Value = "'19.67"
$ws.Range("E45").Value = "  +0.85%  "

# Row 46
$ws.Range("D46").Value = "'0.0509"
$ws.Range("E46").Value = "  -1.46%  "

# Row 47
$ws.Range("D47").Value = "'0.571"
$ws.Range("E47").Value = "  -1.73%  "

# Row 48
$ws.Range("E48").Value = "  -2.81%  "

# Row 49
$ws.Range("E49").Value = "  -0.05%  "

# Row 50
$ws.Range("D50").Value = "'4.67"
$ws.Range("E50").Value = "  +0.06%  "

# Row 51
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").Value = "'1.53"
$ws.Range("E51").Value = "  -2.70%  "

